$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.388616442680359
$ws.Range("B1").Value = 1.664647579193115
$ws.Range("C1").Value = 6.992143154144287
$ws.Range("D1").Value = 1.911772489547729
$ws.Range("E1").Value = 0.8650373220443726
